# Add 2022-Q3 data:
#  1. Duplicate the existing "2022-Q2" sheet (which keeps the same layout as every
#     quarterly "fund detail" sheet) to create a new sheet positioned right before it,
#     rename the new sheet to "2022-Q3", and overwrite its values with the 2022-Q3 figures.
#  2. Insert a new row at the top of the data in "总计" (totals) sheet for the 2022-Q3
#     summary line, pushing the existing quarters down, and fix up the sequential index
#     column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q3" fund-detail sheet
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2, $null)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Columns D,E,F,G hold numeric-looking text (fund size / positions / weights), keep them
# as text just like every other quarter sheet.
$q3.Range("B2:G10").NumberFormat = "@"

# Row 2: 010728 中泰兴诚价值一年持有期混合A
$q3.Range("D2").Value = "6.63"
$q3.Range("E2").Value = "89.26"
$q3.Range("F2").Value = "4.99"
$q3.Range("G2").Value = "0.3308"
$q3.Range("H2").Value = 8

# Row 3: 014772 中泰红利价值一年持有混合
$q3.Range("D3").Value = "5.05"
$q3.Range("E3").Value = "93.55"
$q3.Range("F3").Value = "4.77"
$q3.Range("G3").Value = "0.2409"
$q3.Range("H3").Value = 9

# Row 4: 014771 中泰红利优选一年持有混合
$q3.Range("D4").Value = "5.02"
$q3.Range("E4").Value = "93.69"
$q3.Range("F4").Value = "4.34"
$q3.Range("G4").Value = "0.2179"
$q3.Range("H4").Value = 10

# Row 5: 010729 中泰兴诚价值一年持有期混合C
$q3.Range("D5").Value = "1.27"
$q3.Range("E5").Value = "89.26"
$q3.Range("F5").Value = "4.99"
$q3.Range("G5").Value = "0.0634"
$q3.Range("H5").Value = 8

# Row 6: fund code/name swap to 004317 前海开源沪港深裕鑫灵活配置混合C
$q3.Range("B6").Value = "004317"
$q3.Range("C6").Value = "前海开源沪港深裕鑫灵活配置混合C"
$q3.Range("D6").Value = "1.79"
$q3.Range("E6").Value = "70.17"
$q3.Range("F6").Value = "2.98"
$q3.Range("G6").Value = "0.0533"
$q3.Range("H6").Value = 4

# Row 7: 004316 前海开源沪港深裕鑫灵活配置混合A
$q3.Range("D7").Value = "1.77"
$q3.Range("E7").Value = "70.17"
$q3.Range("F7").Value = "2.98"
$q3.Range("G7").Value = "0.0527"
$q3.Range("H7").Value = 4

# Row 8: fund code/name swap to 004249 安信中国制造混合
$q3.Range("B8").Value = "004249"
$q3.Range("C8").Value = "安信中国制造混合"
$q3.Range("D8").Value = "0.52"
$q3.Range("E8").Value = "89.55"
$q3.Range("F8").Value = "3.78"
$q3.Range("G8").Value = "0.0197"
$q3.Range("H8").Value = 9

# Row 9: 161124 易方达香港恒生综合小型股指数（QDII-LOF）A
$q3.Range("D9").Value = "0.20"
$q3.Range("E9").Value = "91.61"
$q3.Range("F9").Value = "1.35"
$q3.Range("G9").Value = "0.0027"
$q3.Range("H9").Value = 7

# Row 10: 006263 易方达香港恒生综合小型股指数（QDII-LOF）C
$q3.Range("D10").Value = "0.05"
$q3.Range("E10").Value = "91.61"
$q3.Range("F10").Value = "1.35"
$q3.Range("G10").Value = "0.0007"
$q3.Range("H10").Value = 7

# ---------------------------------------------------------------------------
# Step 2: update the "总计" (totals) sheet with the new 2022-Q3 summary row
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Re-apply row3's (the old first data row) formatting onto the freshly inserted
# row2 so the new row keeps the same look (index column style, etc.) instead of
# whatever formatting Insert() guessed at.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 9
$total.Range("D2").Value = 0.98

# Re-sequence the index column (A) for every data row: 0,1,2,...
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7
